# Updates cryptos list values (price/volume/coin swaps) per the
# "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.995.56"
$ws.Range("E2").Value = "  +2.06%  "
$ws.Range("D3").Value = "3.186.69"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  +0.08%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "534.71"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.36%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "145.20"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +3.54%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.21%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.527"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -1.35%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "7.31"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("E11").Value = "  -0.62%  "
$ws.Range("D12").Value = "3.735.24"
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("E13").Value = "  -1.79%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "25.85"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -0.52%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.0000173"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").Value = "60.036.36"
$ws.Range("E16").Value = "  +2.08%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "6.23"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.147.98"
$ws.Range("E18").Value = "  -0.31%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "13.24"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +1.56%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "368.78"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -0.66%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("E23").Value = "  +0.20%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "69.53"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -0.70%  "
$ws.Range("E25").Value = "  +0.29%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "8.67"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +5.92%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("E28").Value = "  +0.07%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "22.27"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +1.10%  "
$ws.Range("E30").Value = "  +0.03%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "5.28"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +1.47%  "
$ws.Range("E33").Value = "  +1.94%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "6.57"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +4.52%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "156.11"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -2.05%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").Value = "2.812.17"
$ws.Range("E37").Value = "  +6.44%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "26.06"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +3.16%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.0703"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +2.97%  "
$ws.Range("E40").Value = "  -0.94%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.0298"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +3.68%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "4.22"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -0.70%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "39.64"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +2.14%  "
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("D46").Value = "3.229.07"
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("E47").Value = "  -0.28%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "6.15"
$cell.Style = "Normal"
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "20.60"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +1.71%  "
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.796"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +5.21%  "
$ws.Range("E51").Value = "  +0.04%  "
